# Applies a cyclic rotation of species-observation data among rows 2-7
# (columns A,B,D,E,F,G,H,Q,R,AC) on the active worksheet, as captured by the
# commit's diff. The content effectively moves:
#   old row 4 -> row 2
#   old row 2 -> row 3
#   old row 3 -> row 7
#   old row 7 -> row 6
#   old row 6 -> row 5
#   old row 5 -> row 4
# (all other columns/rows are left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Snapshot current ("before") values for the rows/columns that participate.
$before = @{}
foreach ($r in 2..7) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $rowData["AC"] = $ws.Range("AC$r").Value2
    $before[$r] = $rowData
}

# Target row gets the snapshot from its source row.
$mapping = @{ 2 = 4; 3 = 2; 4 = 5; 5 = 6; 6 = 7; 7 = 3 }

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $data = $before[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value2 = $data[$c]
    }
    if ($null -eq $data["AC"]) {
        $ws.Range("AC$target").ClearContents()
    } else {
        $ws.Range("AC$target").Value2 = $data["AC"]
    }
}
